$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet had a row-alignment bug: "Docentes responsáveis:" (row 12) had no
# value row under it, which pushed every subsequent label/value pairing off
# by one. Fix it by inserting a new row at 13 for the teacher name, then
# correct the text that is now associated with each label.

$ws.Rows("13:13").Insert()
$ws.Range("A13").Clear()

# Row 13: teacher name now sits under "Docentes responsáveis:" (row 12)
$ws.Range("B13").Value2 = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C13").Value2 = "11079086 - Herlandí de Souza Andrade"

# The freshly-inserted row has no column formatting of its own yet; B13 would
# otherwise fall back to column A's bold style. Pull B13/C13's formatting
# from the row below (already styled correctly for the B/C columns).
$ws.Range("B14:C14").Copy() | Out-Null
$ws.Range("B13:C13").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 14: "Programa resumido:" gets its real summary text (was "Semestral")
$ws.Range("B14").Value2 = "Administração e processo estratégico; Planejamento estratégico, tático e operacional; Diretrizes organizacionais: missão, visão e objetivos; Formulação estratégica; Implantação de estratégia; Controle estratégico; Planejamento de unidades organizacionais; Administração estratégica aplicada."
$ws.Range("C14").Value2 = "Administração e processo estratégico; Planejamento estratégico, tático e operacional; Diretrizes organizacionais: missão, visão e objetivos; Formulação estratégica; Implantação de estratégia; Controle estratégico; Planejamento de unidades organizacionais; Administração estratégica aplicada."

# Row 16: "Programa:" gets its real full syllabus text (was "01/01/2021")
$ws.Range("B16").Value2 = "1. Motivações e Desafios para a estratégia; 2. Conceitos Básicos de Estratégia; 3. Gestão Estratégica; 4. Transformação Estratégica; 5. Análise do Ambiente Externo; 6. Análise da Turbulência e da Vulnerabilidade; 7. Análise do Ambiente Interno; 8. Representação do Portifólio; 9. Estratégia de Balanceamento do Portifólio; 10. Formulação das Estratégias; 11. Capacitação Estratégica; 12. O Plano Estratégico; 13. Metodologia de Planejamento Estratégico; 14. Workshop de Planejamento Estratégico; 15. Implantação da Gestão Estratégica."
$ws.Range("C16").Value2 = "1. Motivações e Desafios para a estratégia; 2. Conceitos Básicos de Estratégia; 3. Gestão Estratégica; 4. Transformação Estratégica; 5. Análise do Ambiente Externo; 6. Análise da Turbulência e da Vulnerabilidade; 7. Análise do Ambiente Interno; 8. Representação do Portifólio; 9. Estratégia de Balanceamento do Portifólio; 10. Formulação das Estratégias; 11. Capacitação Estratégica; 12. O Plano Estratégico; 13. Metodologia de Planejamento Estratégico; 14. Workshop de Planejamento Estratégico; 15. Implantação da Gestão Estratégica."

# Row 19: "Método:" gets its real method text (was the teacher name)
$ws.Range("B19").Value2 = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras."
$ws.Range("C19").Value2 = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras."

# Row 20: "Critério:" gets its real criteria text (was the método text)
$ws.Range("B20").Value2 = "Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas"
$ws.Range("C20").Value2 = "Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas"

# Row 21: "Norma de recuperação:" gets its real recovery-rule text (was the critério text)
$ws.Range("B21").Value2 = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação."
$ws.Range("C21").Value2 = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação."

# Row 22: "Bibliografia:" gets its real bibliography text (was the norma de recuperação text)
$ws.Range("B22").Value2 = "OLIVEIRA, D. P. R. Planejamento Estratégico: Conceitos, Metodologia e Práticas. 34 ed. São Paulo: Atlas, 2018.CHIAVENATO, I; SAPIRO, A. Planejamento Estratégico: Fundamentos e Aplicações. 3 ed. Rio de Janeiro: Campus, 2015.COSTA, E. A. Gestão Estratégica: da empresa que temos para a empresa que queremos. 2 ed. Santo André: Saraiva, 2012.LOBATO, D. M. Estratégia de Empresas. Rio de Janeiro: FGV, 2009.HITT, M A. et al. Administração Estratégica. São Paulo: Pioneira Thomson Learning, 2007.GHEMAWAT, P. A Estratégia e o cenário de Negócios. Porto Alegre: Bookman, 2007.MINTZBERG, H. et al. O Processo da Estratégia. São Paulo: Bookman, 2006.HAMEL, G., PRAHALAD, C.K. Competindo pelo futuro. Rio de Janeiro: Campus, 2005.PORTER, M. Estratégia Competitiva. Rio de janeiro: Campus, 2005.KAPLAN, R. S. Mapas Estratégicos: Balanced Scorecard. Rio de Janeiro: Elsevier, 2004."
$ws.Range("C22").Value2 = "OLIVEIRA, D. P. R. Planejamento Estratégico: Conceitos, Metodologia e Práticas. 34 ed. São Paulo: Atlas, 2018.CHIAVENATO, I; SAPIRO, A. Planejamento Estratégico: Fundamentos e Aplicações. 3 ed. Rio de Janeiro: Campus, 2015.COSTA, E. A. Gestão Estratégica: da empresa que temos para a empresa que queremos. 2 ed. Santo André: Saraiva, 2012.LOBATO, D. M. Estratégia de Empresas. Rio de Janeiro: FGV, 2009.HITT, M A. et al. Administração Estratégica. São Paulo: Pioneira Thomson Learning, 2007.GHEMAWAT, P. A Estratégia e o cenário de Negócios. Porto Alegre: Bookman, 2007.MINTZBERG, H. et al. O Processo da Estratégia. São Paulo: Bookman, 2006.HAMEL, G., PRAHALAD, C.K. Competindo pelo futuro. Rio de Janeiro: Campus, 2005.PORTER, M. Estratégia Competitiva. Rio de janeiro: Campus, 2005.KAPLAN, R. S. Mapas Estratégicos: Balanced Scorecard. Rio de Janeiro: Elsevier, 2004."

# Row 10 ("Objetivos:") had never carried its own Portuguese objectives text
# (it showed the teacher's name); give it the real text too.
$ws.Range("B10").Value2 = "Apresentar aos alunos os fundamentos do planejamento e da gestão estratégica nas organizações, capacitando-os quanto as metodologias existentes, suas etapas e implicações para os resultados organizacionais."
$ws.Range("C10").Value2 = "Apresentar aos alunos os fundamentos do planejamento e da gestão estratégica nas organizações, capacitando-os quanto as metodologias existentes, suas etapas e implicações para os resultados organizacionais."

Write-Output "done"
